$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("D7").Value = 44421
$ws.Range("K7").Value = 'Packham''s Triumph'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 16500
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 917
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44421
$ws.Range("K8").Value = 'Winter Nelis'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 16500
$ws.Range("Q8").Value = '$/bandeja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 917
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44336
$ws.Range("K9").Value = 'Winter Nelis'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("Q9").Value = '$/bandeja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1194
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44371
$ws.Range("K10").Value = 'Packham''s Triumph'
$ws.Range("L10").Value = 'Calibre 90'
$ws.Range("M10").Value = 140
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 17429
$ws.Range("Q10").Value = '$/caja 18 kilos embalada'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 968
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44371
$ws.Range("K11").Value = 'Winter Nelis'
$ws.Range("L11").Value = 'Calibre 80'
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 17500
$ws.Range("Q11").Value = '$/caja 18 kilos embalada'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 972
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44314
$ws.Range("K12").Value = 'Packham''s Triumph'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17500
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 972
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44313
$ws.Range("K13").Value = 'Winter Nelis'
$ws.Range("L13").Value = 'Tercera'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 861
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44292
$ws.Range("K14").Value = 'Packham''s Triumph'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 22000
$ws.Range("O14").Value = 23000
$ws.Range("P14").Value = 22500
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 1250
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("D15").Value = 44292
$ws.Range("K15").Value = 'Winter Nelis'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 22000
$ws.Range("O15").Value = 23000
$ws.Range("P15").Value = 22500
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1250
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("D16").Value = 44398
$ws.Range("K16").Value = 'Packham''s Triumph'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20500
$ws.Range("Q16").Value = '$/caja 20 kilos empedrada'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1025
$ws.Range("T16").Value = 20

# Row 17
$ws.Range("D17").Value = 44398
$ws.Range("K17").Value = 'Winter Nelis'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 20500
$ws.Range("Q17").Value = '$/caja 20 kilos empedrada'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1025
$ws.Range("T17").Value = 20

# Row 18 (new)
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C18").Value = 'Arica y Parinacota'
$ws.Range("D18").Value = 44355
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 'Fruta'
$ws.Range("G18").Value = 100104
$ws.Range("H18").Value = 'Frutos de pepita'
$ws.Range("I18").Value = 100104005
$ws.Range("J18").Value = 'Pera'
$ws.Range("K18").Value = 'Packham''s Triumph'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 17500
$ws.Range("Q18").Value = '$/caja 18 kilos granel'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 972
$ws.Range("T18").Value = 18
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 19 (new)
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C19").Value = 'Arica y Parinacota'
$ws.Range("D19").Value = 44355
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 'Fruta'
$ws.Range("G19").Value = 100104
$ws.Range("H19").Value = 'Frutos de pepita'
$ws.Range("I19").Value = 100104005
$ws.Range("J19").Value = 'Pera'
$ws.Range("K19").Value = 'Winter Nelis'
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 250
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 17500
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 972
$ws.Range("T19").Value = 18
$ws.Range("D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
